$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 246.38461
$ws.Range("I11").Value = 246.38461
$ws.Range("K11").Value = 246.38461
$ws.Range("M11").Value = -106.38461
$ws.Range("H15").Value = 1032.7869
$ws.Range("I15").Value = 1032.7869
$ws.Range("K15").Value = 3098.3607
$ws.Range("M15").Value = -2929.3607
$ws.Range("H40").Value = 71446584
$ws.Range("I40").Value = 2998.5
$ws.Range("K40").Value = 2998.5
$ws.Range("M40").Value = -2823.5
$ws.Range("H43").Value = 3544
$ws.Range("I43").Value = 3099
$ws.Range("J43").Value = 3989
$ws.Range("K43").Value = 3099
$ws.Range("L43").Value = 3989
$ws.Range("M43").Value = -3030
$ws.Range("N43").Value = -4127
$ws.Range("H76").Value = 9480.833000000001
$ws.Range("I76").Value = 16222.75
$ws.Range("J76").Value = 6109.875
$ws.Range("K76").Value = 16222.75
$ws.Range("L76").Value = 6109.875
$ws.Range("M76").Value = -15907.75
$ws.Range("N76").Value = -6739.875
$ws.Range("H79").Value = 9480.833000000001
$ws.Range("I79").Value = 16222.75
$ws.Range("J79").Value = 6109.875
$ws.Range("K79").Value = 16222.75
$ws.Range("L79").Value = 6109.875
$ws.Range("M79").Value = -15130.75
$ws.Range("N79").Value = -8293.875
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H112").Value = 4073.4443
$ws.Range("J112").Value = 4698.8
$ws.Range("L112").Value = 14096.4
$ws.Range("N112").Value = -16312.4
$ws.Range("H129").Value = 4932.7
$ws.Range("I129").Value = 917.8
$ws.Range("K129").Value = 2753.4
$ws.Range("M129").Value = 2246.6
$ws.Range("H135").Value = 7505
$ws.Range("I135").Value = 1574.2
$ws.Range("J135").Value = 17389.666
$ws.Range("K135").Value = 14167.8
$ws.Range("L135").Value = 156506.994
$ws.Range("M135").Value = -11632.8
$ws.Range("N135").Value = -161576.994
$ws.Range("H137").Value = 2406.6785
$ws.Range("I137").Value = 1890.8125
$ws.Range("J137").Value = 3094.5
$ws.Range("K137").Value = 5672.4375
$ws.Range("L137").Value = 9283.5
$ws.Range("M137").Value = -3122.4375
$ws.Range("N137").Value = -14383.5
$ws.Range("H138").Value = 4046.9546
$ws.Range("J138").Value = 4980.1396
$ws.Range("L138").Value = 14940.4188
$ws.Range("N138").Value = -25220.4188

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10031.621
$ws.Range("I32").Value = 9443.223
$ws.Range("K32").Value = 9443.223
$ws.Range("M32").Value = -9156.223
$ws.Range("H61").Value = 5681710
$ws.Range("I61").Value = 6902569
$ws.Range("K61").Value = 6902569
$ws.Range("M61").Value = -6902357
$ws.Range("H122").Value = 9224.75
$ws.Range("I122").Value = 12999.5
$ws.Range("K122").Value = 38998.5
$ws.Range("M122").Value = -36548.5
$ws.Range("H132").Value = 2502711.2
$ws.Range("I132").Value = 2445.6875
$ws.Range("J132").Value = 12503773
$ws.Range("K132").Value = 7337.0625
$ws.Range("L132").Value = 37511319
$ws.Range("M132").Value = -4807.0625
$ws.Range("N132").Value = -37516379
$ws.Range("H136").Value = 5681710
$ws.Range("I136").Value = 6902569
$ws.Range("K136").Value = 20707707
$ws.Range("M136").Value = -20705157

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5003870.5
$ws.Range("I134").Value = 3905.5293
$ws.Range("K134").Value = 11716.5879
$ws.Range("M134").Value = -9181.5879

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34486532
$ws.Range("I31").Value = 71431944
$ws.Range("K31").Value = 71431944
$ws.Range("M31").Value = -71431649
$ws.Range("H34").Value = 34486532
$ws.Range("I34").Value = 71431944
$ws.Range("K34").Value = 71431944
$ws.Range("M34").Value = -71431742
$ws.Range("H58").Value = 2716.524
$ws.Range("J58").Value = 3966.5715
$ws.Range("L58").Value = 3966.5715
$ws.Range("N58").Value = -4372.5715
$ws.Range("H108").Value = 99999
$ws.Range("J108").Value = 99999
$ws.Range("L108").Value = 99999
$ws.Range("N108").Value = -107679
$ws.Range("H115").Value = 54622.25
$ws.Range("J115").Value = 59596.332
$ws.Range("L115").Value = 59596.332
$ws.Range("N115").Value = -61946.332
$ws.Range("H122").Value = 2968.1875
$ws.Range("I122").Value = 2249.2
$ws.Range("K122").Value = 6747.599999999999
$ws.Range("M122").Value = -4297.599999999999
$ws.Range("H132").Value = 1686.7587
$ws.Range("I132").Value = 1813.75
$ws.Range("J132").Value = 1077.2
$ws.Range("K132").Value = 5441.25
$ws.Range("L132").Value = 3231.6
$ws.Range("M132").Value = -2911.25
$ws.Range("N132").Value = -8291.6
$ws.Range("H134").Value = 2212.6775
$ws.Range("I134").Value = 2017.3478
$ws.Range("J134").Value = 2774.25
$ws.Range("K134").Value = 6052.0434
$ws.Range("L134").Value = 8322.75
$ws.Range("M134").Value = -3517.0434
$ws.Range("N134").Value = -13392.75
$ws.Range("H136").Value = 2716.524
$ws.Range("J136").Value = 3966.5715
$ws.Range("L136").Value = 11899.7145
$ws.Range("N136").Value = -16999.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 17780.422
$ws.Range("I122").Value = 37051.555
$ws.Range("J122").Value = 436.4
$ws.Range("K122").Value = 333463.995
$ws.Range("L122").Value = 3927.6
$ws.Range("M122").Value = -331013.995
$ws.Range("N122").Value = -8827.6
$ws.Range("H132").Value = 1720.6666
$ws.Range("I132").Value = 1849.3334
$ws.Range("J132").Value = 1463.3334
$ws.Range("K132").Value = 16644.0006
$ws.Range("L132").Value = 13170.0006
$ws.Range("M132").Value = -14114.0006
$ws.Range("N132").Value = -18230.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 23298
$ws.Range("J54").Value = 23298
$ws.Range("L54").Value = 23298
$ws.Range("N54").Value = -24078
$ws.Range("H70").Value = 14746.36
$ws.Range("I70").Value = 14567.083
$ws.Range("K70").Value = 14567.083
$ws.Range("M70").Value = -14297.083
$ws.Range("H73").Value = 14746.36
$ws.Range("I73").Value = 14567.083
$ws.Range("K73").Value = 14567.083
$ws.Range("M73").Value = -13631.083
$ws.Range("H126").Value = 4125.3657
$ws.Range("I126").Value = 3712.5
$ws.Range("K126").Value = 11137.5
$ws.Range("M126").Value = -8667.5
$ws.Range("H132").Value = 2226081.2
$ws.Range("I132").Value = 3895.4358
$ws.Range("K132").Value = 11686.3074
$ws.Range("M132").Value = -9156.307400000002
$ws.Range("H136").Value = 4995.5
$ws.Range("J136").Value = 4995.5
$ws.Range("L136").Value = 14986.5
$ws.Range("N136").Value = -20086.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 99241
$ws.Range("J6").Value = 99241
$ws.Range("L6").Value = 99241
$ws.Range("N6").Value = -99465
$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H22").Value = 8259491.5
$ws.Range("I22").Value = 18874512
$ws.Range("J22").Value = 3364
$ws.Range("K22").Value = 18874512
$ws.Range("L22").Value = 3364
$ws.Range("M22").Value = -18874217
$ws.Range("N22").Value = -3954
$ws.Range("H27").Value = 8259491.5
$ws.Range("I27").Value = 18874512
$ws.Range("J27").Value = 3364
$ws.Range("K27").Value = 18874512
$ws.Range("L27").Value = 3364
$ws.Range("M27").Value = -18874405
$ws.Range("N27").Value = -3578
$ws.Range("H46").Value = 877.1
$ws.Range("I46").Value = 597.0625
$ws.Range("K46").Value = 597.0625
$ws.Range("M46").Value = -409.0625
$ws.Range("H54").Value = 89999
$ws.Range("J54").Value = 89999
$ws.Range("L54").Value = 89999
$ws.Range("N54").Value = -91287
$ws.Range("H61").Value = 2915.9412
$ws.Range("I61").Value = 2692.6667
$ws.Range("J61").Value = 3451.8
$ws.Range("K61").Value = 2692.6667
$ws.Range("L61").Value = 3451.8
$ws.Range("M61").Value = -2490.6667
$ws.Range("N61").Value = -3855.8
$ws.Range("H113").Value = 2915.9412
$ws.Range("I113").Value = 2692.6667
$ws.Range("J113").Value = 3451.8
$ws.Range("K113").Value = 2692.6667
$ws.Range("L113").Value = 3451.8
$ws.Range("M113").Value = -522.6667000000002
$ws.Range("N113").Value = -7791.8
$ws.Range("H122").Value = 3488.9456
$ws.Range("I122").Value = 3065.796
$ws.Range("J122").Value = 6944.6665
$ws.Range("K122").Value = 9197.387999999999
$ws.Range("L122").Value = 20833.9995
$ws.Range("M122").Value = -6747.387999999999
$ws.Range("N122").Value = -25733.9995
$ws.Range("H132").Value = 3763.738
$ws.Range("I132").Value = 2300.6072
$ws.Range("K132").Value = 6901.821599999999
$ws.Range("M132").Value = -4371.821599999999
$ws.Range("H136").Value = 4469.7085
$ws.Range("I136").Value = 2751.8667
$ws.Range("K136").Value = 8255.6001
$ws.Range("M136").Value = -5705.6001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 47548
$ws.Range("I55").Value = 30064.334
$ws.Range("K55").Value = 30064.334
$ws.Range("M55").Value = -29787.334
$ws.Range("H132").Value = 913482.75
$ws.Range("I132").Value = 4976.375
$ws.Range("J132").Value = 3336166.2
$ws.Range("K132").Value = 14929.125
$ws.Range("L132").Value = 10008498.6
$ws.Range("M132").Value = -12399.125
$ws.Range("N132").Value = -10013558.6
$ws.Range("H136").Value = 258115.6
$ws.Range("I136").Value = 1719.7931
$ws.Range("J136").Value = 1001663.4
$ws.Range("K136").Value = 5159.379300000001
$ws.Range("L136").Value = 3004990.2
$ws.Range("M136").Value = -2609.379300000001
$ws.Range("N136").Value = -3010090.2
